$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.05524200000000001
$ws.Range("H2").Value = 0.165726
$ws.Range("I2").Value = 0.006851115624430413
$ws.Range("J2").Value = 0.006851115624430412
$ws.Range("M2").Value = 0.002913
$ws.Range("N2").Value = 0.008739
$ws.Range("Q2").Value = 0.000160919946
$ws.Range("R2").Value = 0.001448279514
$ws.Range("S2").Value = 0.006851115624430413
$ws.Range("T2").Value = 0.006851115624430412

# Row 3
$ws.Range("I3").Value = 0.1943656285540788
$ws.Range("J3").Value = 0.1943656285540788
$ws.Range("M3").Value = 0.002913
$ws.Range("N3").Value = 0.008739
$ws.Range("Q3").Value = 0.004565286614000001
$ws.Range("R3").Value = 0.04108757952600001
$ws.Range("S3").Value = 0.1943656285540788
$ws.Range("T3").Value = 0.1943656285540788

# Row 4
$ws.Range("G4").Value = 0.3059026666666667
$ws.Range("H4").Value = 0.917708
$ws.Range("I4").Value = 0.03793806413878803
$ws.Range("J4").Value = 0.03793806413878802
$ws.Range("M4").Value = 0.002913
$ws.Range("N4").Value = 0.008739
$ws.Range("Q4").Value = 0.000891094468
$ws.Range("R4").Value = 0.008019850212
$ws.Range("S4").Value = 0.03793806413878803
$ws.Range("T4").Value = 0.03793806413878802

# Row 5
$ws.Range("E5").Value = 1
$ws.Range("F5").Value = 0.3333333333333333
$ws.Range("G5").Value = 0.03313
$ws.Range("H5").Value = 0.09939
$ws.Range("I5").Value = 0.004108784269892103
$ws.Range("J5").Value = 0.004108784269892103
$ws.Range("M5").Value = 0.002913
$ws.Range("N5").Value = 0.008739
$ws.Range("Q5").Value = 0.00009650769000000001
$ws.Range("R5").Value = 0.00086856921
$ws.Range("S5").Value = 0.004108784269892103
$ws.Range("T5").Value = 0.004108784269892103

# Row 6
$ws.Range("G6").Value = 5.575577333333332
$ws.Range("H6").Value = 16.726732
$ws.Range("I6").Value = 0.691483381912676
$ws.Range("J6").Value = 0.691483381912676
$ws.Range("M6").Value = 0.002913
$ws.Range("N6").Value = 0.008739
$ws.Range("Q6").Value = 0.016241656772
$ws.Range("R6").Value = 0.146174910948
$ws.Range("S6").Value = 0.691483381912676
$ws.Range("T6").Value = 0.691483381912676

# Row 7
$ws.Range("G7").Value = 0.526149
$ws.Range("H7").Value = 1.578447
$ws.Range("I7").Value = 0.06525302550013463
$ws.Range("J7").Value = 0.06525302550013462
$ws.Range("M7").Value = 0.002913
$ws.Range("N7").Value = 0.008739
$ws.Range("Q7").Value = 0.001532672037
$ws.Range("R7").Value = 0.013794048333
$ws.Range("S7").Value = 0.06525302550013463
$ws.Range("T7").Value = 0.06525302550013462
